$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 927.14923
$ws.Range("I15").Value = 927.14923
$ws.Range("K15").Value = 2781.44769
$ws.Range("M15").Value = -2612.44769

$ws.Range("H17").Value = 2759.75
$ws.Range("J17").Value = 2759.75
$ws.Range("L17").Value = 8279.25
$ws.Range("N17").Value = -8615.25

$ws.Range("H38").Value = 645.5
$ws.Range("I38").Value = 161.66667
$ws.Range("K38").Value = 485.00001
$ws.Range("M38").Value = -113.00001

$ws.Range("H43").Value = 2238.8333
$ws.Range("I43").Value = 1949.75
$ws.Range("J43").Value = 2321.4285
$ws.Range("K43").Value = 1949.75
$ws.Range("L43").Value = 2321.4285
$ws.Range("M43").Value = -1880.75
$ws.Range("N43").Value = -2459.4285

$ws.Range("H76").Value = 5815.1055
$ws.Range("I76").Value = 5561.875
$ws.Range("K76").Value = 5561.875
$ws.Range("M76").Value = -5246.875

$ws.Range("H79").Value = 5815.1055
$ws.Range("I79").Value = 5561.875
$ws.Range("K79").Value = 5561.875
$ws.Range("M79").Value = -4469.875

$ws.Range("H88").Value = 3987.8667
$ws.Range("I88").Value = 4648.6665
$ws.Range("J88").Value = 3822.6667
$ws.Range("K88").Value = 4648.6665
$ws.Range("L88").Value = 3822.6667
$ws.Range("M88").Value = -4242.6665
$ws.Range("N88").Value = -4634.6667

$ws.Range("H91").Value = 3987.8667
$ws.Range("I91").Value = 4648.6665
$ws.Range("J91").Value = 3822.6667
$ws.Range("K91").Value = 4648.6665
$ws.Range("L91").Value = 3822.6667
$ws.Range("M91").Value = -3244.6665
$ws.Range("N91").Value = -6630.6667

$ws.Range("H98").Value = 1656.8334
$ws.Range("I98").Value = 1338.9286
$ws.Range("K98").Value = 1338.9286
$ws.Range("M98").Value = 159.0714

$ws.Range("H101").Value = 496.42856
$ws.Range("I101").Value = 515
$ws.Range("K101").Value = 1545
$ws.Range("M101").Value = 77

$ws.Range("H103").Value = 511.6
$ws.Range("I103").Value = 339.5
$ws.Range("K103").Value = 1018.5
$ws.Range("M103").Value = -432.5

$ws.Range("H113").Value = 5625.375
$ws.Range("J113").Value = 5498.5
$ws.Range("L113").Value = 5498.5
$ws.Range("N113").Value = -12006.5

$ws.Range("H122").Value = 1656.8334
$ws.Range("I122").Value = 1338.9286
$ws.Range("K122").Value = 4016.7858
$ws.Range("M122").Value = -1566.7858

$ws.Range("H125").Value = 3337
$ws.Range("I125").Value = 3415.6667
$ws.Range("J125").Value = 3258.3333
$ws.Range("K125").Value = 30741.0003
$ws.Range("L125").Value = 29324.9997
$ws.Range("M125").Value = -28281.0003
$ws.Range("N125").Value = -34244.9997

$ws.Range("H132").Value = 65198.62
$ws.Range("I132").Value = 70613.586
$ws.Range("K132").Value = 211840.758
$ws.Range("M132").Value = -209310.758

$ws.Range("H133").Value = 71778.336
$ws.Range("J133").Value = 71778.336
$ws.Range("L133").Value = 71778.336
$ws.Range("N133").Value = -81898.336

$ws.Range("H135").Value = 946.6977
$ws.Range("I135").Value = 933.5238
$ws.Range("K135").Value = 8401.7142
$ws.Range("M135").Value = -5866.7142

$ws.Range("H137").Value = 3320.6
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 3320.6
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 9961.8
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -15061.8

$ws.Range("H138").Value = 1664.1471
$ws.Range("I138").Value = 1077.674
$ws.Range("J138").Value = 2890.4092
$ws.Range("K138").Value = 3233.022
$ws.Range("L138").Value = 8671.2276
$ws.Range("M138").Value = 1906.978
$ws.Range("N138").Value = -18951.2276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1827.1143
$ws.Range("I2").Value = 1972.0385
$ws.Range("K2").Value = 1972.0385
$ws.Range("M2").Value = -1859.0385

$ws.Range("H32").Value = 8476628
$ws.Range("I32").Value = 8622432
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 8622432
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -8622145
$ws.Range("N32").Value = -20574

$ws.Range("H74").Value = 7983.8
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 7983.8
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H88").Value = 1018.0909
$ws.Range("I88").Value = 960.1111
$ws.Range("J88").Value = 1058.2307
$ws.Range("K88").Value = 960.1111
$ws.Range("L88").Value = 1058.2307
$ws.Range("M88").Value = -554.1111
$ws.Range("N88").Value = -1870.2307

$ws.Range("H91").Value = 1018.0909
$ws.Range("I91").Value = 960.1111
$ws.Range("J91").Value = 1058.2307
$ws.Range("K91").Value = 960.1111
$ws.Range("L91").Value = 1058.2307
$ws.Range("M91").Value = 443.8889
$ws.Range("N91").Value = -3866.2307

$ws.Range("H116").Value = 1827.1143
$ws.Range("I116").Value = 1972.0385
$ws.Range("K116").Value = 1972.0385
$ws.Range("M116").Value = 321.9614999999999

$ws.Range("H132").Value = 477225.22
$ws.Range("I132").Value = 543047.56
$ws.Range("J132").Value = 5498.5
$ws.Range("K132").Value = 1629142.68
$ws.Range("L132").Value = 16495.5
$ws.Range("M132").Value = -1626612.68
$ws.Range("N132").Value = -21555.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1827.1143
$ws.Range("I3").Value = 1972.0385
$ws.Range("K3").Value = 1972.0385
$ws.Range("M3").Value = -1858.0385

$ws.Range("H58").Value = 44629.25
$ws.Range("J58").Value = 47985.668
$ws.Range("L58").Value = 47985.668
$ws.Range("N58").Value = -48573.668

$ws.Range("H59").Value = 107500
$ws.Range("J59").Value = 107500
$ws.Range("L59").Value = 107500
$ws.Range("N59").Value = -109194

$ws.Range("H64").Value = 1827.4286
$ws.Range("J64").Value = 1827.4286
$ws.Range("L64").Value = 1827.4286
$ws.Range("N64").Value = -2277.4286

$ws.Range("H67").Value = 1827.4286
$ws.Range("J67").Value = 1827.4286
$ws.Range("L67").Value = 1827.4286
$ws.Range("N67").Value = -3387.4286

$ws.Range("H94").Value = 999.25
$ws.Range("I94").Value = 796.6875
$ws.Range("J94").Value = 2619.75
$ws.Range("K94").Value = 796.6875
$ws.Range("L94").Value = 2619.75
$ws.Range("M94").Value = -345.6875
$ws.Range("N94").Value = -3521.75

$ws.Range("H99").Value = 28523.652
$ws.Range("I99").Value = 28870.158
$ws.Range("K99").Value = 28870.158
$ws.Range("M99").Value = -27372.158

$ws.Range("H105").Value = 2094
$ws.Range("I105").Value = 1731.2
$ws.Range("J105").Value = 3454.5
$ws.Range("K105").Value = 1731.2
$ws.Range("L105").Value = 3454.5
$ws.Range("M105").Value = 15.79999999999995
$ws.Range("N105").Value = -6948.5

$ws.Range("H107").Value = 2817.2222
$ws.Range("J107").Value = 3113
$ws.Range("L107").Value = 3113

$ws.Range("H122").Value = 150000
$ws.Range("J122").Value = 150000
$ws.Range("L122").Value = 150000
$ws.Range("N122").Value = -159800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2455
$ws.Range("I16").Value = 2299.5
$ws.Range("J16").Value = 2766
$ws.Range("K16").Value = 2299.5
$ws.Range("L16").Value = 2766
$ws.Range("M16").Value = -2012.5
$ws.Range("N16").Value = -3340

$ws.Range("H31").Value = 149664.7
$ws.Range("I31").Value = 197726.25
$ws.Range("J31").Value = 53541.625
$ws.Range("K31").Value = 197726.25
$ws.Range("L31").Value = 53541.625
$ws.Range("M31").Value = -197431.25
$ws.Range("N31").Value = -54131.625

$ws.Range("H34").Value = 149664.7
$ws.Range("I34").Value = 197726.25
$ws.Range("J34").Value = 53541.625
$ws.Range("K34").Value = 197726.25
$ws.Range("L34").Value = 53541.625
$ws.Range("M34").Value = -197524.25
$ws.Range("N34").Value = -53945.625

$ws.Range("H41").Value = 4912.625
$ws.Range("I41").Value = 462.2857
$ws.Range("J41").Value = 36065
$ws.Range("K41").Value = 462.2857
$ws.Range("L41").Value = 36065
$ws.Range("M41").Value = -34.28570000000002
$ws.Range("N41").Value = -36921

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H51").Value = 50090
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H58").Value = 282799.16
$ws.Range("I58").Value = 427325.53
$ws.Range("K58").Value = 427325.53
$ws.Range("M58").Value = -427122.53

$ws.Range("H59").Value = 89997
$ws.Range("J59").Value = 89997
$ws.Range("L59").Value = 89997
$ws.Range("N59").Value = -92287

$ws.Range("H60").Value = 17056.4
$ws.Range("J60").Value = 23333.334
$ws.Range("L60").Value = 23333.334
$ws.Range("N60").Value = -24355.334

$ws.Range("H61").Value = 50090
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H99").Value = 5289.273
$ws.Range("I99").Value = 6974.25
$ws.Range("J99").Value = 4326.4287
$ws.Range("K99").Value = 6974.25
$ws.Range("L99").Value = 4326.4287
$ws.Range("M99").Value = -5476.25
$ws.Range("N99").Value = -7322.4287

$ws.Range("H113").Value = 2455
$ws.Range("I113").Value = 2299.5
$ws.Range("J113").Value = 2766
$ws.Range("K113").Value = 2299.5
$ws.Range("L113").Value = 2766
$ws.Range("M113").Value = -129.5
$ws.Range("N113").Value = -7106

$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178

$ws.Range("H126").Value = 5289.273
$ws.Range("I126").Value = 6974.25
$ws.Range("J126").Value = 4326.4287
$ws.Range("K126").Value = 20922.75
$ws.Range("L126").Value = 12979.2861
$ws.Range("M126").Value = -18452.75
$ws.Range("N126").Value = -17919.2861

$ws.Range("H132").Value = 13910278
$ws.Range("I132").Value = 44575
$ws.Range("K132").Value = 133725
$ws.Range("M132").Value = -131195

$ws.Range("H134").Value = 18931.85
$ws.Range("I134").Value = 27182
$ws.Range("J134").Value = 3610.1428
$ws.Range("K134").Value = 81546
$ws.Range("L134").Value = 10830.4284
$ws.Range("M134").Value = -79011
$ws.Range("N134").Value = -15900.4284

$ws.Range("H136").Value = 282799.16
$ws.Range("I136").Value = 427325.53
$ws.Range("K136").Value = 1281976.59
$ws.Range("M136").Value = -1279426.59

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 56
$ws.Range("I12").Value = 51.42857
$ws.Range("K12").Value = 154.28571
$ws.Range("M12").Value = 18.71429000000001

$ws.Range("H60").Value = 1737
$ws.Range("I60").Value = 1032.2727
$ws.Range("J60").Value = 3675
$ws.Range("K60").Value = 3096.8181
$ws.Range("L60").Value = 11025
$ws.Range("M60").Value = -2845.8181
$ws.Range("N60").Value = -11527

$ws.Range("H86").Value = 643.1818
$ws.Range("J86").Value = 625
$ws.Range("L86").Value = 1875
$ws.Range("N86").Value = -4247

$ws.Range("H89").Value = 643.1818
$ws.Range("J89").Value = 625
$ws.Range("L89").Value = 5625
$ws.Range("N89").Value = -17481

$ws.Range("H107").Value = 368
$ws.Range("I107").Value = 233
$ws.Range("K107").Value = 699
$ws.Range("M107").Value = 1221

$ws.Range("H113").Value = 900.05
$ws.Range("I113").Value = 543.6
$ws.Range("J113").Value = 1018.86664
$ws.Range("K113").Value = 1630.8
$ws.Range("L113").Value = 3056.59992
$ws.Range("M113").Value = 539.1999999999998
$ws.Range("N113").Value = -7396.59992

$ws.Range("H131").Value = 7622.4165
$ws.Range("I131").Value = 589.6667
$ws.Range("J131").Value = 9966.667
$ws.Range("K131").Value = 1769.0001
$ws.Range("L131").Value = 29900.001
$ws.Range("M131").Value = 3270.9999
$ws.Range("N131").Value = -39980.001

$ws.Range("H139").Value = 1286.125
$ws.Range("I139").Value = 1286.125
$ws.Range("K139").Value = 3858.375
$ws.Range("M139").Value = 1281.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13336
$ws.Range("J70").Value = 5000
$ws.Range("L70").Value = 5000
$ws.Range("N70").Value = -5540

$ws.Range("H73").Value = 13336
$ws.Range("J73").Value = 5000
$ws.Range("L73").Value = 5000
$ws.Range("N73").Value = -6872

$ws.Range("H97").Value = 2373.9614
$ws.Range("I97").Value = 1110.0526
$ws.Range("K97").Value = 1110.0526
$ws.Range("M97").Value = -614.0526
$ws.Range("N97").Value = -6796.5713

$ws.Range("H102").Value = 2454.3
$ws.Range("I102").Value = 1392.3182
$ws.Range("K102").Value = 1392.3182
$ws.Range("M102").Value = 229.6818000000001

$ws.Range("H113").Value = 3625.1738
$ws.Range("I113").Value = 2693
$ws.Range("J113").Value = 5075.222
$ws.Range("K113").Value = 2693
$ws.Range("L113").Value = 5075.222
$ws.Range("M113").Value = -523
$ws.Range("N113").Value = -9415.222

$ws.Range("H122").Value = 10799.23
$ws.Range("I122").Value = 2762.818
$ws.Range("K122").Value = 8288.454000000002
$ws.Range("M122").Value = -5838.454000000002

$ws.Range("H126").Value = 1519960.1
$ws.Range("J126").Value = 4840.4287
$ws.Range("L126").Value = 14521.2861
$ws.Range("N126").Value = -19461.2861

$ws.Range("H132").Value = 355535.06
$ws.Range("I132").Value = 366263.4
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 1098790.2
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -1096260.2
$ws.Range("N132").Value = -9560

$ws.Range("H133").Value = 99979
$ws.Range("J133").Value = 99979
$ws.Range("L133").Value = 99979
$ws.Range("N133").Value = -110099

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3759.2415
$ws.Range("I16").Value = 4629.2173
$ws.Range("J16").Value = 424.33334
$ws.Range("K16").Value = 4629.2173
$ws.Range("L16").Value = 424.33334
$ws.Range("M16").Value = -4459.2173
$ws.Range("N16").Value = -764.33334

$ws.Range("H61").Value = 4191
$ws.Range("I61").Value = 1958.8
$ws.Range("J61").Value = 6220.273
$ws.Range("K61").Value = 1958.8
$ws.Range("L61").Value = 6220.273
$ws.Range("M61").Value = -1756.8
$ws.Range("N61").Value = -6624.273

$ws.Range("H113").Value = 4191
$ws.Range("I113").Value = 1958.8
$ws.Range("J113").Value = 6220.273
$ws.Range("K113").Value = 1958.8
$ws.Range("L113").Value = 6220.273
$ws.Range("M113").Value = 211.2
$ws.Range("N113").Value = -10560.273

$ws.Range("H132").Value = 1385477.2
$ws.Range("I132").Value = 1385477.2
$ws.Range("K132").Value = 4156431.6
$ws.Range("M132").Value = -4153901.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 16813
$ws.Range("I45").Value = 13000
$ws.Range("K45").Value = 13000
$ws.Range("M45").Value = -12509

$ws.Range("H100").Value = 2299.3684
$ws.Range("I100").Value = 2085.9
$ws.Range("K100").Value = 4171.8
$ws.Range("M100").Value = -3630.8

$ws.Range("H107").Value = 2445.5405
$ws.Range("I107").Value = 1373.4482
$ws.Range("J107").Value = 6331.875
$ws.Range("K107").Value = 4120.3446
$ws.Range("L107").Value = 18995.625
$ws.Range("M107").Value = -2200.3446
$ws.Range("N107").Value = -22835.625

$ws.Range("H113").Value = 1278.0869
$ws.Range("I113").Value = 575.21875
$ws.Range("J113").Value = 2884.6428
$ws.Range("K113").Value = 1725.65625
$ws.Range("L113").Value = 8653.9284
$ws.Range("M113").Value = 444.34375
$ws.Range("N113").Value = -12993.9284

$ws.Range("H126").Value = 2018.7858
$ws.Range("I126").Value = 1784.9
$ws.Range("K126").Value = 5354.700000000001
$ws.Range("M126").Value = -2884.700000000001

$ws.Range("H132").Value = 4109959
$ws.Range("I132").Value = 5751482.5
$ws.Range("J132").Value = 6150.2144
$ws.Range("K132").Value = 17254447.5
$ws.Range("L132").Value = 18450.6432
$ws.Range("M132").Value = -17251917.5
$ws.Range("N132").Value = -23510.6432

$ws.Range("H136").Value = 6916927.5
$ws.Range("I136").Value = 9269470
$ws.Range("J136").Value = 27337.5
$ws.Range("K136").Value = 27808410
$ws.Range("L136").Value = 82012.5
$ws.Range("M136").Value = -27805860
$ws.Range("N136").Value = -87112.5

$ws.Range("H141").Value = 119857
$ws.Range("J141").Value = 119857
$ws.Range("L141").Value = 119857
$ws.Range("N141").Value = -130217

Write-Host "Applied all changes"